$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of an existing data row onto the two new rows so
# the new cells reuse the existing cell style (s="1") instead of Excel
# allocating a brand new style record.
$ws.Range("A9:B9").Copy()
$ws.Range("A10:B11").PasteSpecial(-4122)

$ws.Range("A10").Value = "bitcoin"
$ws.Range("B10").Value = "com.hamxa.shaynachim"
$ws.Range("A11").Value = "passive income ideas"
$ws.Range("B11").Value = "passive.income.nadi.myfirstdrawermenuproject"

$ws.Rows.Item(11).RowHeight = 24

$ws.Range("B11").Select()
